$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 218 (pushes the current rows 218-256 down to 219-257)
$ws.Rows.Item(218).Insert()

# Insert a second new blank row at 246 (post first-insert numbering), which pushes
# the current rows 246-257 down to 247-258. In original numbering this corresponds
# to inserting right before what used to be row 245.
$ws.Rows.Item(246).Insert()

# Populate the brand-new row 218 with its data
$ws.Cells.Item(218, 1).Value = 6
$ws.Cells.Item(218, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(218, 3).Value = "Metropolitana"
$ws.Cells.Item(218, 4).Value = 44748
$ws.Cells.Item(218, 5).Value = 13
$ws.Cells.Item(218, 6).Value = 100112026
$ws.Cells.Item(218, 7).Value = "Haba"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 800
$ws.Cells.Item(218, 11).Value = 16000
$ws.Cells.Item(218, 12).Value = 18000
$ws.Cells.Item(218, 13).Value = 17125
$ws.Cells.Item(218, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(218, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(218, 16).Value = 685
$ws.Cells.Item(218, 17).Value = 25
$ws.Cells.Item(218, 18).Value = "Hortaliza"

# Populate the brand-new row 246 with its data (this is the blank row created above)
$ws.Cells.Item(246, 1).Value = 6
$ws.Cells.Item(246, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(246, 3).Value = "Metropolitana"
$ws.Cells.Item(246, 4).Value = 44747
$ws.Cells.Item(246, 5).Value = 13
$ws.Cells.Item(246, 6).Value = 100112026
$ws.Cells.Item(246, 7).Value = "Haba"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 750
$ws.Cells.Item(246, 11).Value = 18000
$ws.Cells.Item(246, 12).Value = 20000
$ws.Cells.Item(246, 13).Value = 19067
$ws.Cells.Item(246, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(246, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(246, 16).Value = 763
$ws.Cells.Item(246, 17).Value = 25
$ws.Cells.Item(246, 18).Value = "Hortaliza"

Write-Output "done"
